# Add team record (Wins / Losses / Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy the formatting of an existing header cell (A1, which
# carries the bold/centered/bordered header style) onto the three new
# header cells, then set their labels.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows: every player row (2-54) gets the same team record.
$lastRow = 54
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 84  # AD = Wins
    $ws.Cells.Item($r, 31).Value = 78  # AE = Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF = Ties
}
